# Agenda de projeto.xlsx - apply "Criação adas tabs" edit
#
# Summary of the change:
#  - Column F ("Status") values for the "Criar" rows that were still "?"
#    every other data-row (rows 6,8,10,12,14,16,18,20) are updated to "!"
#    (matching rows 4 and 5 which were already "!").
#  - An AutoFilter is applied on the existing A3:G28 filter range:
#       * Column A ("Ação")   -> show only "Criar" and "Teste"
#       * Column F ("Status") -> show only "?"
#    which, combined with the column F edits above, hides rows
#    4-21 and 27-28, leaving rows 1-3 and 22-26 visible.
#  - The sheet view selection moves from B8 to B22, and the view's
#    frozen/scrolled top-left cell (A3) is cleared.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update column F values (rows where "Criar" should now also read "!") ---
foreach ($r in @(6, 8, 10, 12, 14, 16, 18, 20)) {
    $ws.Cells.Item($r, 6).Value = "!"
}

# --- Apply the AutoFilter criteria on the existing A3:G28 range ---
$filterRange = $ws.Range("A3:G28")
$filterRange.AutoFilter(1, @("Criar", "Teste"), 7) | Out-Null
$filterRange.AutoFilter(6, @("?"), 7) | Out-Null

# --- Update the active selection / scroll position ---
$ws.Range("B22").Select() | Out-Null
